# Updates the "Fitness" column (C) values on Sheet1 to reflect the
# corrected best-fitness-so-far series for run_5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2:C3").Value = 12082
$ws.Range("C4").Value = 11808
$ws.Range("C5:C12").Value = 8837
$ws.Range("C13:C23").Value = 8196
$ws.Range("C24:C88").Value = 8146
$ws.Range("C89:C164").Value = 7569
$ws.Range("C165:C252").Value = 7293
